# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 2..5 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1353
$ws1.Range("F3").Value = 2879
$ws1.Range("F4").Value = 3
$ws1.Range("F5").Value = 263

# Sheet "全部类型" - rows 3,4,5,7 in column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1353
$ws4.Range("F4").Value = 2879
$ws4.Range("F5").Value = 3
$ws4.Range("F7").Value = 263
